$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values for rows 2-6 (years 2018-2022 actual/estimate data)
$ws.Range("AA2").Value = 114.45
$ws.Range("AB2").Value = 299.91
$ws.Range("AC2").Value = 1207
$ws.Range("AD2").Value = 13.75
$ws.Range("AE2").Value = 19221
$ws.Range("AF2").Value = 0.86
$ws.Range("AG2").Value = 350
$ws.Range("AH2").Value = 2.11
$ws.Range("AI2").Value = 27.57
$ws.Range("AJ2").Value = 10858393
$ws.Range("D2").Value = 8776
$ws.Range("E2").Value = 710
$ws.Range("F2").Value = 710
$ws.Range("G2").Value = 579
$ws.Range("H2").Value = 407
$ws.Range("I2").Value = 131
$ws.Range("J2").Value = 276
$ws.Range("K2").Value = 9289
$ws.Range("L2").Value = 4958
$ws.Range("M2").Value = 4332
$ws.Range("N2").Value = 1979
$ws.Range("O2").Value = 2353
$ws.Range("P2").Value = 566
$ws.Range("Q2").Value = -55
$ws.Range("R2").Value = -787
$ws.Range("S2").Value = -65
$ws.Range("T2").Value = 191
$ws.Range("U2").Value = -246
$ws.Range("V2").Value = 2305
$ws.Range("W2").Value = 8.09
$ws.Range("X2").Value = 4.63
$ws.Range("Y2").Value = 6.59
$ws.Range("Z2").Value = 4.42
$ws.Range("AA3").Value = 146.21
$ws.Range("AB3").Value = 294.85
$ws.Range("AC3").Value = 1488
$ws.Range("AD3").Value = 10.42
$ws.Range("AE3").Value = 18170
$ws.Range("AF3").Value = 0.85
$ws.Range("AG3").Value = 350
$ws.Range("AH3").Value = 2.26
$ws.Range("AI3").Value = 22.68
$ws.Range("AJ3").Value = 11478870
$ws.Range("D3").Value = 10145
$ws.Range("E3").Value = 454
$ws.Range("F3").Value = 454
$ws.Range("G3").Value = 574
$ws.Range("H3").Value = 419
$ws.Range("I3").Value = 169
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 11425
$ws.Range("L3").Value = 6784
$ws.Range("M3").Value = 4640
$ws.Range("N3").Value = 1984
$ws.Range("O3").Value = 2656
$ws.Range("P3").Value = 597
$ws.Range("Q3").Value = -521
$ws.Range("R3").Value = -435
$ws.Range("S3").Value = 683
$ws.Range("T3").Value = 166
$ws.Range("U3").Value = -687
$ws.Range("V3").Value = 3742
$ws.Range("W3").Value = 4.47
$ws.Range("X3").Value = 4.13
$ws.Range("Y3").Value = 8.51
$ws.Range("Z3").Value = 4.04
$ws.Range("AA4").Value = 147.65
$ws.Range("AB4").Value = 315.9
$ws.Range("AC4").Value = 1117
$ws.Range("AD4").Value = 13.2
$ws.Range("AE4").Value = 18865
$ws.Range("AF4").Value = 0.78
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 2.37
$ws.Range("AI4").Value = 30.02
$ws.Range("AJ4").Value = 11576968
$ws.Range("D4").Value = 8085
$ws.Range("E4").Value = 727
$ws.Range("F4").Value = 533
$ws.Range("G4").Value = 629
$ws.Range("H4").Value = 283
$ws.Range("I4").Value = 128
$ws.Range("J4").Value = 155
$ws.Range("K4").Value = 12635
$ws.Range("L4").Value = 7533
$ws.Range("M4").Value = 5102
$ws.Range("N4").Value = 2078
$ws.Range("O4").Value = 3024
$ws.Range("P4").Value = 602
$ws.Range("Q4").Value = 759
$ws.Range("R4").Value = -764
$ws.Range("S4").Value = 612
$ws.Range("T4").Value = 224
$ws.Range("U4").Value = 535
$ws.Range("V4").Value = 4262
$ws.Range("W4").Value = 8.99
$ws.Range("X4").Value = 3.51
$ws.Range("Y4").Value = 6.32
$ws.Range("Z4").Value = 2.36
$ws.Range("AA5").Value = 156.27
$ws.Range("AB5").Value = 337.52
$ws.Range("AC5").Value = 362
$ws.Range("AD5").Value = 67.63
$ws.Range("AE5").Value = 20013
$ws.Range("AF5").Value = 1.22
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 1.43
$ws.Range("AI5").Value = 100.25
$ws.Range("AJ5").Value = 12373409
$ws.Range("D5").Value = 10995
$ws.Range("E5").Value = 762
$ws.Range("F5").Value = 762
$ws.Range("G5").Value = 570
$ws.Range("H5").Value = 102
$ws.Range("I5").Value = 42
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 14807
$ws.Range("L5").Value = 9029
$ws.Range("M5").Value = 5778
$ws.Range("N5").Value = 2433
$ws.Range("O5").Value = 3345
$ws.Range("P5").Value = 642
$ws.Range("Q5").Value = 1895
$ws.Range("R5").Value = -1378
$ws.Range("S5").Value = 189
$ws.Range("T5").Value = 207
$ws.Range("U5").Value = 1688
$ws.Range("V5").Value = 4356
$ws.Range("W5").Value = 6.93
$ws.Range("X5").Value = 0.93
$ws.Range("Y5").Value = 1.88
$ws.Range("Z5").Value = 0.74
$ws.Range("AA6").Value = 149.7
$ws.Range("AB6").Value = 344.07
$ws.Range("AC6").Value = 270
$ws.Range("AD6").Value = 64.56
$ws.Range("AE6").Value = 20170
$ws.Range("AF6").Value = 0.86
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 2.01
$ws.Range("AI6").Value = 127.92
$ws.Range("AJ6").Value = 12455526
$ws.Range("D6").Value = 12828
$ws.Range("E6").Value = 725
$ws.Range("F6").Value = 725
$ws.Range("G6").Value = 499
$ws.Range("H6").Value = 278
$ws.Range("I6").Value = 34
$ws.Range("K6").Value = 15086
$ws.Range("L6").Value = 9044
$ws.Range("M6").Value = 6042
$ws.Range("N6").Value = 2472
$ws.Range("P6").Value = 646
$ws.Range("Q6").Value = 589
$ws.Range("R6").Value = -452
$ws.Range("S6").Value = -21
$ws.Range("T6").Value = 227
$ws.Range("U6").Value = 363
$ws.Range("V6").Value = 4356
$ws.Range("W6").Value = 5.65
$ws.Range("X6").Value = 2.16
$ws.Range("Y6").Value = 1.37
$ws.Range("Z6").Value = 1.86

# Clear cells for rows 7-9 (2020(E)-2022(E) estimate rows stripped of financial data,
# leaving only index/type/period columns A-C)
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
